$wb = $excel.ActiveWorkbook

# --- Participants sheet: move selection from O1 to M1 (sheet stays inactive) ---
$wsParticipants = $wb.Worksheets.Item("Participants")
$wsParticipants.Range("M1").Select()

# --- Petition sheet: data edits, new row, page setup, selection ---
$wsPetition = $wb.Worksheets.Item("Petition")
$wsPetition.Activate()

# G5 changes from "Active - Submitted to court" (shared string 254) to "n/a" (shared string 12)
$wsPetition.Range("G5").Value = "n/a"

# New row 6, highlighted like row 3-5 (yellow fill == style index 12)
$wsPetition.Range("A6:U6").Interior.Color = 65535
$wsPetition.Range("A6").Value = "testT4245"
$wsPetition.Range("B6").Value = 1
$wsPetition.Range("C6").Value = 4
$wsPetition.Range("D6").Value = "n/a"
$wsPetition.Range("E6").Value = "n/a"
$wsPetition.Range("F6").Value = "n/a"
$wsPetition.Range("G6").Value = "Active - Submitted to court"
$wsPetition.Range("H6").Value = "n/a"
$wsPetition.Range("I6").Value = "n/a"
$wsPetition.Range("J6").Value = "n/a"
$wsPetition.Range("K6").Value = "n/a"
$wsPetition.Range("L6").Value = "n/a"
$wsPetition.Range("M6").Value = "n/a"
$wsPetition.Range("N6").Value = "n/a"
$wsPetition.Range("O6").Value = "n/a"
$wsPetition.Range("P6").Value = "n/a"
$wsPetition.Range("Q6").Value = "n/a"
$wsPetition.Range("R6").Value = "n/a"
$wsPetition.Range("S6").Value = "n/a"
$wsPetition.Range("T6").Value = "n/a"
$wsPetition.Range("U6").Value = "n/a"

# Touch page setup so a <pageSetup/> element is emitted (portrait, as in the diff)
$wsPetition.PageSetup.Orientation = 1

# Final selection on the Petition sheet
$wsPetition.Range("U6").Select()

Write-Host "done"
